$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first 4 data rows (rows 2 through 5), shifting the remaining
# data rows up and shrinking the used range from A1:E42 to A1:E38.
$ws.Range("A2:E5").EntireRow.Delete()
